$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add week 35 (AL1) and week 36 (AM1) ---
$ws.Range("AL1").Value = "35"
$ws.Range("AM1").Value = "36"

# --- Row 28: corrected counts for existing weeks + new R28 cell ---
$ws.Range("O28").Value = 1
$ws.Range("P28").Value = 0
$ws.Range("R28").Value = 0
$ws.Range("S28").Value = 1
$ws.Range("T28").Value = 8

# --- New week 35 (AL) / week 36 (AM) data columns ---
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0
$ws.Range("AL4").Value = 0
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0
$ws.Range("AL6").Value = 2
$ws.Range("AM6").Value = 2
$ws.Range("AL7").Value = 0
$ws.Range("AM7").Value = 0
$ws.Range("AL8").Value = 0
$ws.Range("AM8").Value = 0
$ws.Range("AL9").Value = 0
$ws.Range("AM9").Value = 0
$ws.Range("AL10").Value = 0
$ws.Range("AM10").Value = 0
$ws.Range("AM11").Value = 0
$ws.Range("AM12").Value = 0
$ws.Range("AL13").Value = 0
$ws.Range("AM13").Value = 0
$ws.Range("AL14").Value = 0
$ws.Range("AM14").Value = 0
$ws.Range("AM15").Value = 0
$ws.Range("AL16").Value = 0
$ws.Range("AM16").Value = 0
$ws.Range("AL17").Value = 0
$ws.Range("AM17").Value = 0
$ws.Range("AM18").Value = 0
$ws.Range("AM19").Value = 0
$ws.Range("AL22").Value = 0
$ws.Range("AM22").Value = 0
$ws.Range("AL23").Value = 0
$ws.Range("AM23").Value = 0
$ws.Range("AL25").Value = 0
$ws.Range("AM25").Value = 0
$ws.Range("AL26").Value = 0
$ws.Range("AM27").Value = 0
$ws.Range("AL28").Value = 0
$ws.Range("AM28").Value = 0
$ws.Range("AL29").Value = 1
$ws.Range("AM29").Value = 1
$ws.Range("AL30").Value = 14
$ws.Range("AM30").Value = 5
$ws.Range("AL31").Value = 0
$ws.Range("AL35").Value = 2
$ws.Range("AL36").Value = 0
$ws.Range("AM36").Value = 0
$ws.Range("AL37").Value = 0
$ws.Range("AM37").Value = 0
$ws.Range("AL38").Value = 0
$ws.Range("AM38").Value = 0
$ws.Range("AL40").Value = 0
$ws.Range("AM40").Value = 0
$ws.Range("AL41").Value = 0
$ws.Range("AM41").Value = 0
$ws.Range("AL42").Value = 0
$ws.Range("AM42").Value = 0
$ws.Range("AL43").Value = 0
$ws.Range("AL44").Value = 0
$ws.Range("AM44").Value = 0
$ws.Range("AL45").Value = 0
$ws.Range("AM45").Value = 0
$ws.Range("AL46").Value = 0
$ws.Range("AM46").Value = 0
$ws.Range("AL47").Value = 0
$ws.Range("AM47").Value = 0
$ws.Range("AL48").Value = 0
$ws.Range("AM48").Value = 0
$ws.Range("AL49").Value = 0
$ws.Range("AM49").Value = 0
$ws.Range("AL50").Value = 0
$ws.Range("AM50").Value = 0
$ws.Range("AL51").Value = 0
$ws.Range("AM51").Value = 0
$ws.Range("AL52").Value = 0
$ws.Range("AL53").Value = 0
$ws.Range("AM53").Value = 0
$ws.Range("AL54").Value = 0
$ws.Range("AM54").Value = 0
$ws.Range("AL55").Value = 0
$ws.Range("AM55").Value = 0
$ws.Range("AL56").Value = 0
$ws.Range("AM56").Value = 0
$ws.Range("AL57").Value = 0
$ws.Range("AM57").Value = 0
$ws.Range("AL58").Value = 0
$ws.Range("AM58").Value = 0
